$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the duplicate value in C2, then merge B2:C2
$ws.Range("C2").Value = $null
$ws.Range("B2:C2").Merge()
$ws.Range("B2:C2").HorizontalAlignment = -4131  # xlLeft

# Update the selection to the merged cell
$ws.Range("B2:C2").Select()
